$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" -- the b.md file has now been handed off for
# zh-cn and de-de, so its status flips from "Handed back: in sync with en-US"
# to "Ready for handoff", a new (later) handoff xliff + datetime is recorded,
# and an error message is surfaced because the handback file that originally
# came in is stale relative to the new source.
# ---------------------------------------------------------------------------

$readyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f07472630ee848f8607847920ea268589f9f4225/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a35d8b70c442e872d315cec1581af0328b0e47f/e2e/b.md."

# --- Overview sheet: b.md row (row 3) ---------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value2 = $readyForHandoff
$overview.Range("F3").Value2 = $readyForHandoff
$overview.Range("G3").Value2 = "2016-08-23 20:39:31"

# --- zh-cn sheet: b.md row (row 3) -------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value2 = $readyForHandoff
# Leading "'" forces text (not an auto-converted Boolean) for the literal
# "True"/"False" strings used in this column, matching the source data type.
$zhcn.Range("F3").Value2 = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value2 = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value2 = "2016-08-23 20:39:26"
$zhcn.Range("P3").Value2 = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet: b.md row (row 3) -------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value2 = $readyForHandoff
$dede.Range("F3").Value2 = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value2 = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value2 = "2016-08-23 20:39:31"
$dede.Range("P3").Value2 = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.15
